$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.163080333333333
$ws.Range("H2").Value = 3.489241
$ws.Range("I2").Value = 0.004254373353458465
$ws.Range("J2").Value = 0.004254373353458466
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6385226666666667
$ws.Range("N2").Value = 1.915568
$ws.Range("O2").Value = 0.291843422078211
$ws.Range("P2").Value = 0.291843422078211
$ws.Range("Q2").Value = 0.7426531559875555
$ws.Range("R2").Value = 6.683878403887999
$ws.Range("S2").Value = 0.001241610878271673
$ws.Range("T2").Value = 0.001241610878271673

# Row 3
$ws.Range("G3").Value = 1.163080333333333
$ws.Range("H3").Value = 3.489241
$ws.Range("I3").Value = 0.004254373353458465
$ws.Range("J3").Value = 0.004254373353458466
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6443276666666667
$ws.Range("N3").Value = 1.932983
$ws.Range("O3").Value = 0.2944966576696867
$ws.Range("P3").Value = 0.2944966576696867
$ws.Range("Q3").Value = 0.7494048373225556
$ws.Range("R3").Value = 6.744643535903
$ws.Range("S3").Value = 0.001252898733072495
$ws.Range("T3").Value = 0.001252898733072495

# Row 4
$ws.Range("G4").Value = 1.163080333333333
$ws.Range("H4").Value = 3.489241
$ws.Range("I4").Value = 0.004254373353458465
$ws.Range("J4").Value = 0.004254373353458466
$ws.Range("M4").Value = 0.9050443333333335
$ws.Range("N4").Value = 2.715133
$ws.Range("O4").Value = 0.4136599202521024
$ws.Range("P4").Value = 0.4136599202521024
$ws.Range("Q4").Value = 1.052639264894778
$ws.Range("R4").Value = 9.473753384052999
$ws.Range("S4").Value = 0.001759863742114298
$ws.Range("T4").Value = 0.001759863742114298

# Row 5
$ws.Range("G5").Value = 126.8932496666667
$ws.Range("H5").Value = 380.679749
$ws.Range("I5").Value = 0.464156468511879
$ws.Range("J5").Value = 0.464156468511879
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6385226666666667
$ws.Range("N5").Value = 1.915568
$ws.Range("O5").Value = 0.291843422078211
$ws.Range("P5").Value = 0.291843422078211
$ws.Range("Q5").Value = 81.02421615915912
$ws.Range("R5").Value = 729.217945432432
$ws.Range("S5").Value = 0.1354610121502441
$ws.Range("T5").Value = 0.1354610121502441

# Row 6
$ws.Range("G6").Value = 126.8932496666667
$ws.Range("H6").Value = 380.679749
$ws.Range("I6").Value = 0.464156468511879
$ws.Range("J6").Value = 0.464156468511879
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6443276666666667
$ws.Range("N6").Value = 1.932983
$ws.Range("O6").Value = 0.2944966576696867
$ws.Range("P6").Value = 0.2944966576696867
$ws.Range("Q6").Value = 81.76083147347413
$ws.Range("R6").Value = 735.847483261267
$ws.Range("S6").Value = 0.1366925286125135
$ws.Range("T6").Value = 0.1366925286125136

# Row 7
$ws.Range("G7").Value = 126.8932496666667
$ws.Range("H7").Value = 380.679749
$ws.Range("I7").Value = 0.464156468511879
$ws.Range("J7").Value = 0.464156468511879
$ws.Range("M7").Value = 0.9050443333333335
$ws.Range("N7").Value = 2.715133
$ws.Range("O7").Value = 0.4136599202521024
$ws.Range("P7").Value = 0.4136599202521024
$ws.Range("Q7").Value = 114.8440165490686
$ws.Range("R7").Value = 1033.596148941617
$ws.Range("S7").Value = 0.1920029277491213
$ws.Range("T7").Value = 0.1920029277491213

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.00603
$ws.Range("H8").Value = 0.01809
$ws.Range("I8").Value = 0.00002205683527278959
$ws.Range("J8").Value = 0.0000220568352727896
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6385226666666667
$ws.Range("N8").Value = 1.915568
$ws.Range("O8").Value = 0.291843422078211
$ws.Range("P8").Value = 0.291843422078211
$ws.Range("Q8").Value = 0.00385029168
$ws.Range("R8").Value = 0.03465262512
$ws.Range("S8").Value = 0.000006437142286226305
$ws.Range("T8").Value = 0.000006437142286226306

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.00603
$ws.Range("H9").Value = 0.01809
$ws.Range("I9").Value = 0.00002205683527278959
$ws.Range("J9").Value = 0.0000220568352727896
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.6443276666666667
$ws.Range("N9").Value = 1.932983
$ws.Range("O9").Value = 0.2944966576696867
$ws.Range("P9").Value = 0.2944966576696867
$ws.Range("Q9").Value = 0.00388529583
$ws.Range("R9").Value = 0.03496766247
$ws.Range("S9").Value = 0.000006495664266607389
$ws.Range("T9").Value = 0.00000649566426660739

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.00603
$ws.Range("H10").Value = 0.01809
$ws.Range("I10").Value = 0.00002205683527278959
$ws.Range("J10").Value = 0.0000220568352727896
$ws.Range("M10").Value = 0.9050443333333335
$ws.Range("N10").Value = 2.715133
$ws.Range("O10").Value = 0.4136599202521024
$ws.Range("P10").Value = 0.4136599202521024
$ws.Range("Q10").Value = 0.00545741733
$ws.Range("R10").Value = 0.04911675597
$ws.Range("S10").Value = 0.000009124028719955902
$ws.Range("T10").Value = 0.000009124028719955904

# Row 11
$ws.Range("G11").Value = 32.943215
$ws.Range("H11").Value = 98.829645
$ws.Range("I11").Value = 0.1205013377464496
$ws.Range("J11").Value = 0.1205013377464497
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6385226666666667
$ws.Range("N11").Value = 1.915568
$ws.Range("O11").Value = 0.291843422078211
$ws.Range("P11").Value = 0.291843422078211
$ws.Range("Q11").Value = 21.03498949037333
$ws.Range("R11").Value = 189.31490541336
$ws.Range("S11").Value = 0.03516752277292615
$ws.Range("T11").Value = 0.03516752277292616

# Row 12
$ws.Range("G12").Value = 32.943215
$ws.Range("H12").Value = 98.829645
$ws.Range("I12").Value = 0.1205013377464496
$ws.Range("J12").Value = 0.1205013377464497
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.6443276666666667
$ws.Range("N12").Value = 1.932983
$ws.Range("O12").Value = 0.2944966576696867
$ws.Range("P12").Value = 0.2944966576696867
$ws.Range("Q12").Value = 21.22622485344834
$ws.Range("R12").Value = 191.036023681035
$ws.Range("S12").Value = 0.03548724121105548
$ws.Range("T12").Value = 0.03548724121105549

# Row 13
$ws.Range("G13").Value = 32.943215
$ws.Range("H13").Value = 98.829645
$ws.Range("I13").Value = 0.1205013377464496
$ws.Range("J13").Value = 0.1205013377464497
$ws.Range("M13").Value = 0.9050443333333335
$ws.Range("N13").Value = 2.715133
$ws.Range("O13").Value = 0.4136599202521024
$ws.Range("P13").Value = 0.4136599202521024
$ws.Range("Q13").Value = 29.81507005753167
$ws.Range("R13").Value = 268.335630517785
$ws.Range("S13").Value = 0.04984657376246801
$ws.Range("T13").Value = 0.04984657376246802

# Row 14
$ws.Range("G14").Value = 112.3790663333333
$ws.Range("H14").Value = 337.137199
$ws.Range("I14").Value = 0.4110657635529401
$ws.Range("J14").Value = 0.4110657635529401
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.6385226666666667
$ws.Range("N14").Value = 1.915568
$ws.Range("O14").Value = 0.291843422078211
$ws.Range("P14").Value = 0.291843422078211
$ws.Range("Q14").Value = 71.75658111267023
$ws.Range("R14").Value = 645.809230014032
$ws.Range("S14").Value = 0.1199668391344828
$ws.Range("T14").Value = 0.1199668391344828

# Row 15
$ws.Range("G15").Value = 112.3790663333333
$ws.Range("H15").Value = 337.137199
$ws.Range("I15").Value = 0.4110657635529401
$ws.Range("J15").Value = 0.4110657635529401
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.6443276666666667
$ws.Range("N15").Value = 1.932983
$ws.Range("O15").Value = 0.2944966576696867
$ws.Range("P15").Value = 0.2944966576696867
$ws.Range("Q15").Value = 72.40894159273523
$ws.Range("R15").Value = 651.6804743346171
$ws.Range("S15").Value = 0.1210574934487786
$ws.Range("T15").Value = 0.1210574934487786

# Row 16
$ws.Range("G16").Value = 112.3790663333333
$ws.Range("H16").Value = 337.137199
$ws.Range("I16").Value = 0.4110657635529401
$ws.Range("J16").Value = 0.4110657635529401
$ws.Range("M16").Value = 0.9050443333333335
$ws.Range("N16").Value = 2.715133
$ws.Range("O16").Value = 0.4136599202521024
$ws.Range("P16").Value = 0.4136599202521024
$ws.Range("Q16").Value = 101.7080371702741
$ws.Range("R16").Value = 915.3723345324671
$ws.Range("S16").Value = 0.1700414309696787
$ws.Range("T16").Value = 0.1700414309696788
